$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the replacement values look like plain numbers (no significance
# stars), so force those specific cells to text format first to keep them
# stored as text (matching the shared-string table), not numeric values.
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"

$ws.Range("B2").Value = "-0.29***"
$ws.Range("B3").Value = "-1.02*"
$ws.Range("B4").Value = "0.06***"
$ws.Range("B5").Value = "0.01***"

$ws.Range("C2").Value = "-0.02***"
$ws.Range("C3").Value = "-0.42***"
$ws.Range("C4").Value = "-0.0*"
$ws.Range("C5").Value = "0.0***"

$ws.Range("D2").Value = "0.26*"
$ws.Range("D3").Value = "1.65"
$ws.Range("D4").Value = "0.32***"
$ws.Range("D5").Value = "-0.01*"

$ws.Range("E2").Value = "-3.7*"
$ws.Range("E3").Value = "-1.49"
$ws.Range("E4").Value = "0.53"
$ws.Range("E5").Value = "0.1"
